$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ================= Sheet1 (Overview) =================
$ws1.Range("A2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws1.Range("A3").Value = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
$ws1.Range("A4").Value = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
$ws1.Range("B4").Value = "Handed back: in sync with en-US"
$ws1.Range("C4").Value = "Handed back: in sync with en-US"
$ws1.Range("A2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md", "", "", "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/a8b4ec37-3611-47ba-989c-b0307579875a.md", "", "", "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/.localization-config", "", "", ".localization-config")

# ================= Sheet2 (zh-cn) =================
$ws2.Range("A2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws2.Range("C2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-10 01:17:38"
$ws2.Range("E2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws2.Range("F2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-10 01:18:28"
$ws2.Range("A3").Value = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
$ws2.Range("A4").Value = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
$ws2.Range("B4").Value = "Handed back: in sync with en-US"
$ws2.Range("C4").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-10 01:15:40"
$ws2.Range("E4").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$ws2.Range("F4").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-03-10 01:16:20"
$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/05d194503093fbd283883c17b37f0ce33af4a026/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/89e320d4efd0fa33fdffb07b59d8f3cbe9a71e3d/e2e/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/58695a001e38e9c8fea401c1267b3ab253f57f3c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md", "", "", "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/05d194503093fbd283883c17b37f0ce33af4a026/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/89e320d4efd0fa33fdffb07b59d8f3cbe9a71e3d/e2e/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/58695a001e38e9c8fea401c1267b3ab253f57f3c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/a8b4ec37-3611-47ba-989c-b0307579875a.md", "", "", "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1686fa95eba8cffc77ede3e5b39b1f915fc39ded/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/89e320d4efd0fa33fdffb07b59d8f3cbe9a71e3d/e2e/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md")
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/58695a001e38e9c8fea401c1267b3ab253f57f3c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/.localization-config", "", "", ".localization-config")

# ================= Sheet3 (de-de) =================
$ws3.Range("A2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws3.Range("C2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-10 01:17:44"
$ws3.Range("E2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws3.Range("F2").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-10 01:18:44"
$ws3.Range("A3").Value = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
$ws3.Range("A4").Value = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
$ws3.Range("B4").Value = "Handed back: in sync with en-US"
$ws3.Range("C4").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-10 01:15:46"
$ws3.Range("E4").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$ws3.Range("F4").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
$ws3.Range("G4").Value = "2016-03-10 01:16:38"
$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/070c20bb87cfd1cb0a7a19263edcd05aa8609905/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/057263aaa1e31e74c31c73cb65fb1b8f7e1e7f4d/e2e/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/13e0e46baf40a588db024937325dfe2e8f0df42c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf", "", "", "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md", "", "", "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/070c20bb87cfd1cb0a7a19263edcd05aa8609905/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/057263aaa1e31e74c31c73cb65fb1b8f7e1e7f4d/e2e/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/13e0e46baf40a588db024937325dfe2e8f0df42c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/e2e/a8b4ec37-3611-47ba-989c-b0307579875a.md", "", "", "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd2118a4d2cf546559dd14980639279ac04cca6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/057263aaa1e31e74c31c73cb65fb1b8f7e1e7f4d/e2e/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md")
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/13e0e46baf40a588db024937325dfe2e8f0df42c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf", "", "", "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4bdd6e94016651f204e04aa97ec6a5a0e3415442/.localization-config", "", "", ".localization-config")
